$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text looks like a plain number need an explicit
# Text number format first, otherwise Excel auto-converts the typed
# string into a numeric value (losing things like trailing zeros).

$ws.Range("D2").Value = '23.783.51'
$ws.Range("E2").Value = '  -3.57%  '

$ws.Range("D3").Value = '1.612.14'
$ws.Range("E3").Value = '  -4.02%  '

$ws.Range("E4").Value = '  +0.40%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.66'
$ws.Range("E5").Value = '  -2.25%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3890'
$ws.Range("E7").Value = '  -0.88%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3808'
$ws.Range("E8").Value = '  -3.85%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.000'
$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.354'
$ws.Range("E10").Value = '  -3.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.79'
$ws.Range("E11").Value = '  -4.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08417'
$ws.Range("E12").Value = '  -3.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.82'
$ws.Range("E13").Value = '  -6.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.004'
$ws.Range("E14").Value = '  -4.54%  '

$ws.Range("E15").Value = '  -3.90%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.473'
$ws.Range("E16").Value = '  -3.46%  '

$ws.Range("D17").Value = '1.600.93'
$ws.Range("E17").Value = '  -4.48%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.97'
$ws.Range("E18").Value = '  -0.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06906'
$ws.Range("E19").Value = '  -1.91%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.00'
$ws.Range("E20").Value = '  -6.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.788'
$ws.Range("E21").Value = '  -4.13%  '

$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.40'
$ws.Range("E23").Value = '  -4.12%  '

$ws.Range("D24").Value = '23.844.04'
$ws.Range("E24").Value = '  -3.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.437'
$ws.Range("E25").Value = '  +3.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.840'
$ws.Range("E26").Value = '  +2.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.10'
$ws.Range("E27").Value = '  -4.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.31'
$ws.Range("E28").Value = '  -2.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '138.87'
$ws.Range("E29").Value = '  -5.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.246'
$ws.Range("E30").Value = '  -10.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.821'
$ws.Range("E31").Value = '  -6.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.488'
$ws.Range("E32").Value = '  -1.13%  '

$ws.Range("D33").Value = '1.783.53'
$ws.Range("E33").Value = '  -3.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08065'
$ws.Range("E34").Value = '  -3.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9690'
$ws.Range("E35").Value = '  -2.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02864'
$ws.Range("E36").Value = '  -7.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.591'
$ws.Range("E37").Value = '  -5.56%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2649'
$ws.Range("E38").Value = '  -5.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09161'
$ws.Range("E39").Value = '  -3.81%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.31'
$ws.Range("E40").Value = '  -0.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.62'
$ws.Range("E41").Value = '  -0.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.424'
$ws.Range("E42").Value = '  -6.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7437'
$ws.Range("E43").Value = '  -6.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.81'
$ws.Range("E44").Value = '  -5.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6814'
$ws.Range("E45").Value = '  -4.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.438'
$ws.Range("E46").Value = '  -4.76%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.056'
$ws.Range("E47").Value = '  -2.83%  '

$ws.Range("E48").Value = '  +0.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08218'
$ws.Range("E49").Value = '  -5.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.72'
$ws.Range("E50").Value = '  -3.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.212'
$ws.Range("E51").Value = '  -8.84%  '
